# Weekly fruit/vegetable update: add two new price records (rows 228-229)
# for "Ají" (Inferno, Primera / Segunda) at Terminal La Palmera de La Serena,
# pushing the existing rows 228-253 down to 230-255.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 228, shifting everything
# below (old rows 228-253) down to rows 230-255.
$ws.Range("A228:R229").EntireRow.Insert()

# New row 228: Ají / Inferno / Primera
$ws.Range("A228").Value = 8
$ws.Range("B228").Value = "Terminal La Palmera de La Serena"
$ws.Range("C228").Value = "Coquimbo"
$ws.Range("D228").Value = 44769
$ws.Range("E228").Value = 4
$ws.Range("F228").Value = 100112021
$ws.Range("G228").Value = "Ají"
$ws.Range("H228").Value = "Inferno"
$ws.Range("I228").Value = "Primera"
$ws.Range("J228").Value = 520
$ws.Range("K228").Value = 14000
$ws.Range("L228").Value = 15000
$ws.Range("M228").Value = 14500
$ws.Range("N228").Value = "$/caja 12 kilos"
$ws.Range("O228").Value = "Región de Arica y Parinacota"
$ws.Range("P228").Value = 1208
$ws.Range("Q228").Value = 12
$ws.Range("R228").Value = "Hortaliza"

# New row 229: Ají / Inferno / Segunda
$ws.Range("A229").Value = 8
$ws.Range("B229").Value = "Terminal La Palmera de La Serena"
$ws.Range("C229").Value = "Coquimbo"
$ws.Range("D229").Value = 44769
$ws.Range("E229").Value = 4
$ws.Range("F229").Value = 100112021
$ws.Range("G229").Value = "Ají"
$ws.Range("H229").Value = "Inferno"
$ws.Range("I229").Value = "Segunda"
$ws.Range("J229").Value = 300
$ws.Range("K229").Value = 9000
$ws.Range("L229").Value = 10000
$ws.Range("M229").Value = 9500
$ws.Range("N229").Value = "$/caja 12 kilos"
$ws.Range("O229").Value = "Región de Arica y Parinacota"
$ws.Range("P229").Value = 792
$ws.Range("Q229").Value = 12
$ws.Range("R229").Value = "Hortaliza"
